$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.187.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "'3.421.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'412.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "'128.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("E7").Value = "  -2.33%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.728"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("D10").Value = "'0.139"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").Value = "'42.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'0.0000219"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "'9.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "'3.956.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "'20.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.00%  "
$ws.Range("D17").Value = "'3.420.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "'12.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.30%  "
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").Value = "'62.157.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "'472.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.07%  "
$ws.Range("D22").Value = "'91.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("D24").Value = "'13.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'3.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "'9.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.66%  "
$ws.Range("D27").Value = "'33.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").Value = "'7.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'11.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.94%  "
$ws.Range("D32").Value = "'0.165"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.39%  "
$ws.Range("E33").Value = "  -3.46%  "
$ws.Range("D34").Value = "'40.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.81%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "'57.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.09%  "
$ws.Range("D37").Value = "'0.0486"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.77%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").Value = "'3.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.34%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'146.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.28%  "
$ws.Range("B43").Value = "LidoDAOToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D43").Value = "'3.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'2.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.85%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'4.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.20%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'2.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.05%  "
$ws.Range("E47").Value = "  +18.36%  "
$ws.Range("D48").Value = "'16.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("D49").Value = "'0.0₃0541"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +27.29%  "
$ws.Range("D50").Value = "'22.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "'112.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.08%  "
